$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 149 (shifts existing rows 149-232 down to 150-233)
$ws.Rows.Item(149).Insert()

# Populate the new row 149 with data
$ws.Cells.Item(149, 1).Value = 5
$ws.Cells.Item(149, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(149, 3).Value = "Maule"
$ws.Cells.Item(149, 4).Value = 44518
$ws.Cells.Item(149, 5).Value = 7
$ws.Cells.Item(149, 6).Value = 100112032
$ws.Cells.Item(149, 7).Value = "Zapallo italiano"
$ws.Cells.Item(149, 8).Value = "Sin especificar"
$ws.Cells.Item(149, 9).Value = "Primera"
$ws.Cells.Item(149, 10).Value = 400
$ws.Cells.Item(149, 11).Value = 6000
$ws.Cells.Item(149, 12).Value = 6000
$ws.Cells.Item(149, 13).Value = 6000
$ws.Cells.Item(149, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(149, 15).Value = "Región del Maule"
$ws.Cells.Item(149, 16).Value = 100
$ws.Cells.Item(149, 17).Value = 60
$ws.Cells.Item(149, 18).Value = "Hortaliza"

Write-Output "done"
